# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" column (E16:E47) previously listed periods in
# descending order (2103 .. 1808). The database refresh re-sorts this
# block in ascending order (1808 .. 2103), which is the "part 1" of the
# new account-statement periods. The "Valor Mora" column (F16:F47) stays
# tied to its period: every period keeps F = 33019 except the most
# recent period "2103", which keeps F = 25315 - only now that pairing
# rides along with period 2103's new row (47) instead of row 16.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Ascending list of periods that should now populate E16:E47 (top to bottom)
$periods = @(
    "1808","1809","1810","1811","1812",
    "1901","1902","1903","1904","1905","1906","1907","1908","1909","1910","1911","1912",
    "2001","2002","2003","2004","2005","2006","2007","2008","2009","2010","2011","2012",
    "2101","2102","2103"
)

$startRow = 16
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = $startRow + $i
    $period = $periods[$i]

    $ws.Cells.Item($row, 5).Value = $period

    if ($period -eq "2103") {
        $ws.Cells.Item($row, 6).Value = 25315
    } else {
        $ws.Cells.Item($row, 6).Value = 33019
    }
}
